$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $range = $d.Content
    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Text not found: $oldText"
    }
}

# 1. First reflection paragraph - "R/ Al observar la Carta Gantt..."
Replace-Text "R/ Al observar la Carta Gantt y compararlo con nuestros avances puedo decir que vamos algo atrasados " "R/ Al observar la Carta Gantt y compararlo con nuestros avances puedo decir que en la semana 7 íbamos algo atrasados "

Replace-Text " Sprint ya que vamos en el Sprint 1 y según nuestr" " Sprint, ya que íbamos en el Sprint 1 y según nuestr"

Replace-Text "a Carta Gantt deberíamos ir en el 2, esto fue mas que nada porque nos dimos la semana de fiestas patrias, pero en la semana 8 y 9 tuvimos un muy buen avance en relaciona a la programación del proyecto y la documentación de este" "a Carta Gantt deberíamos ir en el 2, esto fue más que nada porque nos dimos la semana de fiestas patrias, pero en la semana 8, 9 y 10, tuvimos un muy buen avance en relaciona a la programación del proyecto y la documentación de este, logrando avanzar rápido y quedar casi al día con las actividades correspondientes a nuestro backlog"

# 2. "R/Según mi trabajo en el grupo..."
Replace-Text "R/Según mi trabajo en el grupo creo que e avanzado bien y de forma eficiente ya que e ayudado en todos los documentos y la programación del proyecto que estamos realizando, siento que quizás podría mejorar en relación a la preparación y organización de las fechas de entregas ya que en ese tema soy algo disperso." "R/Según mi trabajo en el grupo creo que he avanzado bien y de forma eficiente ya que he ayudado en todos los documentos y la programación del proyecto que estamos realizando, siento que quizás podría mejorar en relación a la preparación y organización en relación a que somos más llevados a nuestro ritmo que realizar actividades todos los días como una rutina."

# 3. "De momento no creo tener inquietudes..."
Replace-Text " De momento no creo tener inquietudes, ya que todas las que se me han presentado hasta la fecha las e consultado con mi equipo de capston " " De momento no creo tener inquietudes, ya que todas las que se me han presentado hasta la fecha las he consultado con mi equipo de capston "

# 4. "En mi consideración personal creo que..."
Replace-Text " En mi consideración personal creo que las actividades las debemos escoger cada uno ya que los 3 sabemos que somos un grupo y que debemos ser responsables, además de que al escoger cada uno la actividad se puede trabajar de una forma mas cómoda y que le guste a cada uno." " En mi consideración personal creo que las actividades las debemos escoger cada uno ya que los 3 sabemos que somos un grupo y que debemos ser responsables, además de que al escoger cada uno la actividad se puede trabajar de una forma más cómoda y que le guste a cada uno."

# 5. "En todo el desarrollo del proyecto..."
Replace-Text "En todo el desarrollo del proyecto se nos presentan nuevas actividades pero como grupo nos hemos comunicado bien y dividido estas parte para desarrollar el avance el este proyecto durante la semana 8 y 9" "En todo el desarrollo del proyecto se nos presentan nuevas actividades, pero como grupo nos hemos comunicado bien y dividido estas parte para desarrollar el avance el este proyecto durante la semana 8, 9 y 10"

# 6. "nos a costado decidirnos..."
Replace-Text "nos a costado decidirnos en algún momento en como empezar o quien empieza con que parte pero esto se resolvió durante el mismo día." "nos ha costado decidirnos en algún momento en como empezar o quien empieza con qué parte, pero esto se resolvió durante el mismo día."

Write-Output "Done"
